{"js": "// Replace the date line and the multiplication-table answers with the\n// new values from the commit. Every \"before\" string is unique within\n// the document, so a plain text search + replace is sufficient and\n// keeps existing run formatting (font, size) untouched.\nconst replacements = [\n  [\"2024-11-08 Friday\", \"2024-11-09 Saturday\"],\n  [\"60\u00d760=3600\", \"97\u00d761=5917\"],\n  [\"88\u00d723=2024\", \"79\u00d742=3318\"],\n  [\"76\u00d729=2204\", \"35\u00d727=945\"],\n  [\"35\u00d733=1155\", \"97\u00d729=2813\"],\n  [\"71\u00d792=6532\", \"80\u00d733=2640\"],\n  [\"84\u00d795=7980\", \"94\u00d763=5922\"],\n  [\"27\u00d723=621\", \"70\u00d725=1750\"],\n  [\"74\u00d726=1924\", \"57\u00d772=4104\"],\n  [\"17\u00d792=1564\", \"34\u00d790=3060\"],\n  [\"87\u00d747=4089\", \"96\u00d768=6528\"],\n  [\"88\u00d717=1496\", \"73\u00d754=3942\"],\n  [\"41\u00d752=2132\", \"22\u00d751=1122\"],\n  [\"80\u00d778=6240\", \"51\u00d742=2142\"],\n  [\"67\u00d723=1541\", \"62\u00d781=5022\"],\n  [\"73\u00d788=6424\", \"83\u00d730=2490\"],\n  [\"38\u00d795=3610\", \"39\u00d746=1794\"],\n  [\"32\u00d739=1248\", \"92\u00d723=2116\"],\n  [\"25\u00d744=1100\", \"68\u00d779=5372\"],\n  [\"23\u00d719=437\", \"88\u00d754=4752\"],\n  [\"94\u00d751=4794\", \"26\u00d731=806\"],\n  [\"46\u00d774=3404\", \"64\u00d718=1152\"],\n  [\"61\u00d764=3904\", \"33\u00d765=2145\"],\n  [\"28\u00d763=1764\", \"53\u00d713=689\"],\n  [\"38\u00d759=2242\", \"80\u00d725=2000\"],\n  [\"28\u00d717=476\", \"86\u00d799=8514\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + before);\n  }\n\n  for (const item of results.items) {\n    item.insertText(after, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Each \"before\" string is unique within the document, so a single\n# wdReplaceAll Find/Execute per pair performs exactly one substitution\n# while preserving the existing run formatting (font, size, etc.).\n$replacements = @(\n    @(\"2024-11-08 Friday\", \"2024-11-09 Saturday\"),\n    @(\"60\u00d760=3600\", \"97\u00d761=5917\"),\n    @(\"88\u00d723=2024\", \"79\u00d742=3318\"),\n    @(\"76\u00d729=2204\", \"35\u00d727=945\"),\n    @(\"35\u00d733=1155\", \"97\u00d729=2813\"),\n    @(\"71\u00d792=6532\", \"80\u00d733=2640\"),\n    @(\"84\u00d795=7980\", \"94\u00d763=5922\"),\n    @(\"27\u00d723=621\", \"70\u00d725=1750\"),\n    @(\"74\u00d726=1924\", \"57\u00d772=4104\"),\n    @(\"17\u00d792=1564\", \"34\u00d790=3060\"),\n    @(\"87\u00d747=4089\", \"96\u00d768=6528\"),\n    @(\"88\u00d717=1496\", \"73\u00d754=3942\"),\n    @(\"41\u00d752=2132\", \"22\u00d751=1122\"),\n    @(\"80\u00d778=6240\", \"51\u00d742=2142\"),\n    @(\"67\u00d723=1541\", \"62\u00d781=5022\"),\n    @(\"73\u00d788=6424\", \"83\u00d730=2490\"),\n    @(\"38\u00d795=3610\", \"39\u00d746=1794\"),\n    @(\"32\u00d739=1248\", \"92\u00d723=2116\"),\n    @(\"25\u00d744=1100\", \"68\u00d779=5372\"),\n    @(\"23\u00d719=437\", \"88\u00d754=4752\"),\n    @(\"94\u00d751=4794\", \"26\u00d731=806\"),\n    @(\"46\u00d774=3404\", \"64\u00d718=1152\"),\n    @(\"61\u00d764=3904\", \"33\u00d765=2145\"),\n    @(\"28\u00d763=1764\", \"53\u00d713=689\"),\n    @(\"38\u00d759=2242\", \"80\u00d725=2000\"),\n    @(\"28\u00d717=476\", \"86\u00d799=8514\"),\n)\n\nforeach ($pair in $replacements) {\n    $before = $pair[0]\n    $after = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($before, $false, $true, $false, $false, $false, $true, 0, $false, $after, 2)\n    if (-not $found) {\n        throw \"Could not find text to replace: $before\"\n    }\n}\n"}
